# Applies the "repositories.pptx" doc edit:
#  1. Slide 1, shape "TextBox 86" (id=87): the "5. develop" box -
#     the git commit line gets a "-s" (sign-off) flag added, and the
#     textbox is widened to fit the new text.
#  2. Every "datetimeFigureOut" date placeholder (slide master + all
#     11 custom layouts) gets its cached text bumped from 12/8/2021
#     to 7/15/2022.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) "git commit -a -m " -> "git commit -s -a -m " + widen textbox
# ---------------------------------------------------------------
$slide = $p.Slides.Item(1)
$slideShapes = $slide.Shapes

$devBox = $null
for ($i = 1; $i -le $slideShapes.Count; $i++) {
    $candidate = $slideShapes.Item($i)
    if ($candidate.Id -eq 87) {
        $devBox = $candidate
    }
}

$para = $devBox.TextFrame.TextRange.Paragraphs(4, 1)
$run1 = $para.Runs(1, 1)
$run1.Text = "git commit –s –a –m "

# Widen the textbox to match the longer line (cx 1415772 -> 1553630 EMU).
# PowerPoint COM units are points (914400 EMU/in = 12700 EMU/pt); use a
# value whose point->EMU floor lands exactly on 1553630.
$devBox.Width = 122.3331

# ---------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" field text everywhere
# ---------------------------------------------------------------
$newDate = "7/15/2022"

$master = $p.SlideMaster
$masterShapes = $master.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $shp = $masterShapes.Item($j)
    if ($shp.PlaceholderFormat.Type -eq 16) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $shp = $layoutShapes.Item($j)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
